# se separa campo nombre en apellido y nombre de personal
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$newRow = 56
$ws.Cells.Item($newRow, 1).Value = "anular remito pendiente de fc"
$ws.Cells.Item($newRow, 2).Value = "no comenzado"

$ws.Range("C53").Select()
